$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 124
$ws.Range("H124").Value = 28388
$ws.Range("J124").Value = 28388
$ws.Range("L124").Value = 28388
$ws.Range("N124").Value = -38208
# Row 138
$ws.Range("H138").Value = 2337.63
$ws.Range("I138").Value = 940.0909
$ws.Range("J138").Value = 3435.6965
$ws.Range("K138").Value = 2820.2727
$ws.Range("L138").Value = 10307.0895
$ws.Range("M138").Value = 2319.7273
$ws.Range("N138").Value = -20587.0895
# Row 140
$ws.Range("H140").Value = 48176.832
$ws.Range("J140").Value = 48176.832
$ws.Range("L140").Value = 48176.832
$ws.Range("N140").Value = -58536.832
# Row 141
$ws.Range("H141").Value = 2168.9404
$ws.Range("I141").Value = 1546.3334
$ws.Range("J141").Value = 2463.8596
$ws.Range("K141").Value = 4639.0002
$ws.Range("L141").Value = 7391.578799999999
$ws.Range("M141").Value = 540.9997999999996
$ws.Range("N141").Value = -17751.5788

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2863.111
$ws.Range("I2").Value = 1745.3572
$ws.Range("J2").Value = 6775.25
$ws.Range("K2").Value = 1745.3572
$ws.Range("L2").Value = 6775.25
$ws.Range("M2").Value = -1632.3572
$ws.Range("N2").Value = -7001.25
# Row 44
$ws.Range("H44").Value = 32032.666
$ws.Range("J44").Value = 32032.666
$ws.Range("L44").Value = 32032.666
$ws.Range("N44").Value = -33008.666
# Row 55
$ws.Range("H55").Value = 39582
$ws.Range("J55").Value = 39582
$ws.Range("L55").Value = 39582
$ws.Range("N55").Value = -40212
# Row 61
$ws.Range("H61").Value = 2052.76
$ws.Range("I61").Value = 1892.1305
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 1892.1305
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -1680.1305
$ws.Range("N61").Value = -4324
# Row 74
$ws.Range("H74").Value = 25330.139
$ws.Range("I74").Value = 1280.2285
$ws.Range("K74").Value = 1280.2285
$ws.Range("M74").Value = -406.2284999999999
# Row 77
$ws.Range("H77").Value = 25330.139
$ws.Range("I77").Value = 1280.2285
$ws.Range("K77").Value = 6401.1425
$ws.Range("M77").Value = -2033.1425
# Row 80
$ws.Range("H80").Value = 32984
$ws.Range("J80").Value = 38955
$ws.Range("L80").Value = 38955
$ws.Range("N80").Value = -40951
# Row 83
$ws.Range("H83").Value = 32984
$ws.Range("J83").Value = 38955
$ws.Range("L83").Value = 116865
$ws.Range("N83").Value = -126849
# Row 102
$ws.Range("H102").Value = 71430410
$ws.Range("I102").Value = 1872.1111
$ws.Range("K102").Value = 1872.1111
$ws.Range("M102").Value = -250.1111000000001
# Row 116
$ws.Range("H116").Value = 2863.111
$ws.Range("I116").Value = 1745.3572
$ws.Range("J116").Value = 6775.25
$ws.Range("K116").Value = 1745.3572
$ws.Range("L116").Value = 6775.25
$ws.Range("M116").Value = 548.6428000000001
$ws.Range("N116").Value = -11363.25
# Row 136
$ws.Range("H136").Value = 2052.76
$ws.Range("I136").Value = 1892.1305
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 5676.3915
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -3126.3915
$ws.Range("N136").Value = -16800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2863.111
$ws.Range("I3").Value = 1745.3572
$ws.Range("J3").Value = 6775.25
$ws.Range("K3").Value = 1745.3572
$ws.Range("L3").Value = 6775.25
$ws.Range("M3").Value = -1631.3572
$ws.Range("N3").Value = -7003.25
# Row 35
$ws.Range("H35").Value = 32821.145
$ws.Range("J35").Value = 32821.145
$ws.Range("L35").Value = 32821.145
$ws.Range("N35").Value = -33441.145
# Row 82
$ws.Range("H82").Value = 24977.525
$ws.Range("I82").Value = 9410.875
$ws.Range("J82").Value = 36298.727
$ws.Range("K82").Value = 9410.875
$ws.Range("L82").Value = 36298.727
$ws.Range("M82").Value = -9027.875
$ws.Range("N82").Value = -37064.727
# Row 85
$ws.Range("H85").Value = 24977.525
$ws.Range("I85").Value = 9410.875
$ws.Range("J85").Value = 36298.727
$ws.Range("K85").Value = 9410.875
$ws.Range("L85").Value = 36298.727
$ws.Range("M85").Value = -8084.875
$ws.Range("N85").Value = -38950.727
# Row 130
$ws.Range("H130").Value = 51962.5
$ws.Range("J130").Value = 51962.5
$ws.Range("L130").Value = 51962.5
$ws.Range("N130").Value = -62002.5
# Row 134
$ws.Range("H134").Value = 1682.3043
$ws.Range("I134").Value = 1320.421
$ws.Range("J134").Value = 3401.25
$ws.Range("K134").Value = 3961.263
$ws.Range("L134").Value = 10203.75
$ws.Range("M134").Value = -1426.263
$ws.Range("N134").Value = -15273.75
# Row 135
$ws.Range("H135").Value = 73096.664
$ws.Range("J135").Value = 73096.664
$ws.Range("L135").Value = 73096.664
$ws.Range("N135").Value = -83236.664

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 20729.75
$ws.Range("J41").Value = 24306.666
$ws.Range("L41").Value = 24306.666
$ws.Range("N41").Value = -25162.666
# Row 51
$ws.Range("H51").Value = 9279.333000000001
$ws.Range("J51").Value = 9279.333000000001
$ws.Range("L51").Value = 9279.333000000001
$ws.Range("N51").Value = -10751.333
# Row 61
$ws.Range("H61").Value = 9279.333000000001
$ws.Range("J61").Value = 9279.333000000001
$ws.Range("L61").Value = 9279.333000000001
$ws.Range("N61").Value = -9975.333000000001
# Row 68
$ws.Range("H68").Value = 16999.5
$ws.Range("J68").Value = 16999.5
$ws.Range("L68").Value = 16999.5
$ws.Range("N68").Value = -18497.5
# Row 71
$ws.Range("H71").Value = 16999.5
$ws.Range("J71").Value = 16999.5
$ws.Range("L71").Value = 50998.5
$ws.Range("N71").Value = -58486.5
# Row 131
$ws.Range("H131").Value = 41435.332
$ws.Range("J131").Value = 41435.332
$ws.Range("L131").Value = 41435.332
$ws.Range("N131").Value = -51515.332

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 604.5161000000001
$ws.Range("I113").Value = 583.5833
$ws.Range("J113").Value = 617.7368
$ws.Range("K113").Value = 1750.7499
$ws.Range("L113").Value = 1853.2104
$ws.Range("M113").Value = 419.2501
$ws.Range("N113").Value = -6193.2104
# Row 139
$ws.Range("H139").Value = 21704.96
$ws.Range("I139").Value = 1613.125
$ws.Range("J139").Value = 55543.844
$ws.Range("K139").Value = 4839.375
$ws.Range("L139").Value = 166631.532
$ws.Range("M139").Value = 300.625
$ws.Range("N139").Value = -176911.532

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 11989.643
$ws.Range("J93").Value = 11989.643
$ws.Range("L93").Value = 11989.643
$ws.Range("N93").Value = -15733.643
# Row 123
$ws.Range("H123").Value = 14965.923
$ws.Range("J123").Value = 14965.923
$ws.Range("L123").Value = 14965.923
$ws.Range("N123").Value = -19865.923

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3250
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 4750
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 4750
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -5022
# Row 61
$ws.Range("H61").Value = 2164.4546
$ws.Range("I61").Value = 1500.8125
$ws.Range("J61").Value = 3934.1667
$ws.Range("K61").Value = 1500.8125
$ws.Range("L61").Value = 3934.1667
$ws.Range("M61").Value = -1298.8125
$ws.Range("N61").Value = -4338.1667
# Row 109
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774
# Row 113
$ws.Range("H113").Value = 2164.4546
$ws.Range("I113").Value = 1500.8125
$ws.Range("J113").Value = 3934.1667
$ws.Range("K113").Value = 1500.8125
$ws.Range("L113").Value = 3934.1667
$ws.Range("M113").Value = 669.1875
$ws.Range("N113").Value = -8274.1667
# Row 132
$ws.Range("H132").Value = 4512
$ws.Range("I132").Value = 4588.737
$ws.Range("J132").Value = 4329.75
$ws.Range("K132").Value = 13766.211
$ws.Range("L132").Value = 12989.25
$ws.Range("M132").Value = -11236.211
$ws.Range("N132").Value = -18049.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 48972.332
$ws.Range("J135").Value = 48972.332
$ws.Range("L135").Value = 48972.332
$ws.Range("N135").Value = -59112.332
